$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.534.22'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.47%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.833.32'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.62%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.90'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.27%  '

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.07%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4266'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.39%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3647'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.46%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07254'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.90%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8625'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.64%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.63'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.04%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.797.22'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -2.09%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.409'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.18%  '

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.09%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06950'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.02%  '

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.07%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '80.44'

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008889'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.69%  '

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.02%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.39'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.24%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '27.445.46'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.63%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.146'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +3.31%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.052.42'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.44%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.994'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.25%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '154.67'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.67%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.81'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +1.54%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.121'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.55%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '114.36'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -4.28%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.819'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.72%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08845'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.44%  '

$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7459'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.83%  '

$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = 'HuobiToken'
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.983'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.64%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.531'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.18%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.128'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.60%  '

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.03%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.090'

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05302'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -2.52%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01933'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.02%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.802'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.84%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5070'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.15%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.475'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.45%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.302'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.83%  '

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.43%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.06479'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.94%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '105.21'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.63%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4663'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.39%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.0000'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.11%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.615'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.14%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '63.42'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.64%  '
